$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ============================================================
# Table 2 (task list) gets its new bottom row created first, so
# that its text lands in the shared-string table ahead of the
# other new strings introduced below.
# ============================================================
$ws.Cells.Item(18,1).Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Cells.Item(18,2).Copy()
$ws.Range("B20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(20,1).Value2 = "козметични промени"
$ws.Cells.Item(20,2).Value2 = "менюта"

# ============================================================
# Table 1 (GitHub handles / participants), currently rows 1-8.
# A new header row is inserted above it, pushing the 8 data
# rows down to rows 2-9.
# ============================================================

# Shift the existing A1:C8 block down by one row (values + the
# formatting already sitting in the destination rows, which is
# a leftover mirror of the same table one row down).
$ws.Range("A2:C9").Value2 = $ws.Range("A1:C8").Value2

# Write the new header into row 1
$ws.Cells.Item(1,1).Value2 = "GitHub"
$ws.Cells.Item(1,2).Value2 = "Име"
$ws.Cells.Item(1,3).Value2 = "TAG"

# The old "ок" tags (originally rows 1,5,7,8) should not survive
# the shift - remove them completely (now at rows 2,6,8,9)
$ws.Cells.Item(2,3).Clear()
$ws.Cells.Item(6,3).Clear()
$ws.Cells.Item(8,3).Clear()
$ws.Cells.Item(9,3).Clear()

# Rows 5 and 7 (previously "ок" rows, unstyled) now hold TAG
# values that should carry the same red-font styling as the
# rest of column C - copy it over from a cell that already has it
$ws.Cells.Item(4,3).Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ============================================================
# Table 2 (task list), currently rows 12-18.
# A new header row is inserted above it (still at row 12),
# pushing the data rows down to rows 13-19, and one new task
# row is appended at the end (row 20).
# ============================================================

$ws.Range("A13:B19").Value2 = $ws.Range("A12:B18").Value2

# Make sure the shifted-down row 19 keeps the highlighted style
# used throughout this table (row 20 already has it, see above)
$ws.Cells.Item(12,1).Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Cells.Item(12,2).Copy()
$ws.Range("B19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write the new header into row 12 (no special formatting, like
# the header row of table 1)
$ws.Cells.Item(12,1).Value2 = "задача"
$ws.Cells.Item(12,2).Value2 = "TAG"
$ws.Range("A12:B12").ClearFormats()

# ============================================================
# Final selection, as recorded in the saved workbook
# ============================================================
$ws.Range("C11").Select()
